$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing AutoFilter so it can be cleanly re-applied over the
# expanded header range once the new columns are in place.
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Insert a new column "Sluiten" just before "Actiepunten Projectleider"
# (originally column R), as confirmation from project leaders that
# something can be closed.
$ws.Columns("R").Insert()
$ws.Range("R1").Value = "Sluiten"

# Insert a new column "Warning" right after "Informatie" (which shifted
# from column U to column V after the previous insert).
$ws.Columns("W").Insert()
$ws.Range("W1").Value = "Warning"

# Give the two brand-new columns a sensible width matching their
# neighbours (Q's width for the "Sluiten" column, V's width for the
# "Warning" column).
$ws.Columns("R").ColumnWidth = 17.830729166666668
$ws.Columns("W").ColumnWidth = 29.944010416666668

# Re-apply the AutoFilter across the new full header range (it now spans
# through column Z instead of X).
$null = $ws.Range("A1:Z1").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new
# AutoFilter range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Overzicht!_FilterDatabase") {
        $n.RefersTo = "=Overzicht!`$A`$1:`$Z`$1"
    }
}

# Update the active selection to reflect where the user ended up working.
$ws.Activate()
$null = $ws.Range("U15").Select()
